{"js": "// 1) Remove the leftover \"_GoBack\" bookmark (Word strips this automatically\n//    on save in a real session; here we delete it explicitly).\nconst goBackExists = context.document.bookmarks.exists(\"_GoBack\");\nawait context.sync();\nif (goBackExists.value) {\n  context.document.deleteBookmark(\"_GoBack\");\n  await context.sync();\n}\n\n// 2) Fix the computed discharge time: \"0.217 ns\" -> \"0.679 ns\"\n//    (per commit message: \"Fixed part e&f based on new current/cap calc\").\nconst results = context.document.body.search(\" = 2.1737e-10 seconds = 0.217 ns\", {\n  matchCase: true,\n  matchWholeWord: false\n});\nresults.load(\"text\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  const target = results.items[0];\n  target.insertText(\" = 2.1737e-10 seconds = 0.679 ns\", Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) Remove the leftover \"_GoBack\" bookmark.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n# 2) Fix the computed discharge time: \"0.217 ns\" -> \"0.679 ns\"\n#    (per commit message: \"Fixed part e&f based on new current/cap calc\").\n$find = $d.Content.Find\n$find.Text = \" = 2.1737e-10 seconds = 0.217 ns\"\n$find.Replacement.Text = \" = 2.1737e-10 seconds = 0.679 ns\"\n$find.Execute([ref]\" = 2.1737e-10 seconds = 0.217 ns\", [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$true, [ref]1, [ref]$false, [ref]\" = 2.1737e-10 seconds = 0.679 ns\", [ref]2)\n"}
